# Fix Training Data Issue (#48)
# The "Date" column (BF) on the Team-Data sheet held an incorrectly
# formatted/shifted date string ("2-16-2007-08") for every data row.
# Correct it to the proper ISO date representation ("2008-02-16"),
# keeping the value a plain text string (matching the original inlineStr
# cell type) rather than letting Excel auto-convert it to a date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 31
$dateColumn = 58   # column BF

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateColumn)
    $cell.NumberFormat = "@"
    $cell.Value = "2008-02-16"
}
